# Rename the header row from the old "_old"/"_new" suffix convention to the
# new "_FV2404"/"_FV2410" (format-version) suffix convention, turn the data
# range into a proper Excel Table ("Table1"), and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base header labels (without the old/new suffix) in column order.
$baseHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J (1-10): "<name>_old" -> "<name>_FV2404"
for ($i = 0; $i -lt $baseHeaders.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($baseHeaders[$i])_FV2404"
}

# Column K (11) holds "diff" and is left untouched.

# Columns L-U (12-21): "<name>_new" -> "<name>_FV2410"
for ($i = 0; $i -lt $baseHeaders.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($baseHeaders[$i])_FV2410"
}

# Turn A1:U55 into an Excel Table ("Table1") with a header row, autofilter
# and banded rows - this also wires up the worksheet <tableParts>,
# xl/tables/table1.xml, the relationship and the content-type override.
$tableRange = $ws.Range("A1:U55")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# Freeze the header row (split under row 1, pane state = frozen).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
